# A new team ("Pepe Leal FC") is inserted into the sorted time_id table on
# Sheet1. It takes over the time_id that "Texas Club 2026" used to have
# (1273719) and "Texas Club 2026" is reassigned a new, larger time_id
# (1326835). Every row from the old "Texas Club 2026" row onward therefore
# shifts down by one row, and a brand-new last row is added for the team
# that used to be last.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 19 down to 20, 18 down to 19, ... 6 down to 7, working from the
# bottom up so no data is overwritten before it is copied. (Using Value2 for
# reads since Value's getter is not reliable in this environment; Value is
# fine for writes.)
for ($r = 19; $r -ge 6; $r--) {
    $dest = $r + 1
    $ws.Cells.Item($dest, 1).Value = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($dest, 2).Value = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($dest, 3).Value = $ws.Cells.Item($r, 3).Value2
}

# Fill the freed-up row 6 with the new team's data.
$ws.Cells.Item(6, 1).Value = 1273719
$ws.Cells.Item(6, 2).Value = "Pepe Leal FC"
$ws.Cells.Item(6, 3).Value = 0

# "Texas Club 2026" (now row 7) receives its new time_id.
$ws.Cells.Item(7, 1).Value = 1326835

# Make sure the new cells in column A carry the same style as the rest of
# the column (bordered/centered), and update the sheet's used range.
$ws.Cells.Item(19, 1).Copy()
$ws.Cells.Item(20, 1).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(6, 1).Copy()
$ws.Cells.Item(6, 1).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
